$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3741.2222
$ws.Range("J17").Value = 3889.4
$ws.Range("L17").Value = 11668.2
$ws.Range("N17").Value = -12004.2
$ws.Range("H40").Value = 1500
$ws.Range("J40").Value = 1500
$ws.Range("L40").Value = 1500
$ws.Range("N40").Value = -1850
$ws.Range("H64").Value = 3633.3333
$ws.Range("I64").Value = 4266.6665
$ws.Range("J64").Value = 3000
$ws.Range("K64").Value = 4266.6665
$ws.Range("L64").Value = 3000
$ws.Range("M64").Value = -4018.6665
$ws.Range("N64").Value = -3496
$ws.Range("H67").Value = 3633.3333
$ws.Range("I67").Value = 4266.6665
$ws.Range("J67").Value = 3000
$ws.Range("K67").Value = 4266.6665
$ws.Range("L67").Value = 3000
$ws.Range("M67").Value = -3408.6665
$ws.Range("N67").Value = -4716
$ws.Range("H74").Value = 106170
$ws.Range("I74").Value = 7712.5
$ws.Range("K74").Value = 7712.5
$ws.Range("M74").Value = -6776.5
$ws.Range("H76").Value = 4550.75
$ws.Range("I76").Value = 4567.6665
$ws.Range("K76").Value = 4567.6665
$ws.Range("M76").Value = -4252.6665
$ws.Range("H77").Value = 106170
$ws.Range("I77").Value = 7712.5
$ws.Range("K77").Value = 38562.5
$ws.Range("M77").Value = -33882.5
$ws.Range("H79").Value = 4550.75
$ws.Range("I79").Value = 4567.6665
$ws.Range("K79").Value = 4567.6665
$ws.Range("M79").Value = -3475.6665
$ws.Range("H80").Value = 1003.7273
$ws.Range("I80").Value = 1027.4445
$ws.Range("K80").Value = 3082.3335
$ws.Range("M80").Value = -2084.3335
$ws.Range("H83").Value = 1003.7273
$ws.Range("I83").Value = 1027.4445
$ws.Range("K83").Value = 9247.0005
$ws.Range("M83").Value = -4255.0005

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H15").Value = 4988.8887
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 4988.8887
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 4988.8887
$ws.Range("M15").ClearContents()
$ws.Range("N15").Value = -5688.8887
$ws.Range("H63").Value = 1969.7273
$ws.Range("J63").Value = 2779.6
$ws.Range("L63").Value = 2779.6
$ws.Range("N63").Value = -4151.6
$ws.Range("H66").Value = 1969.7273
$ws.Range("J66").Value = 2779.6
$ws.Range("L66").Value = 13898
$ws.Range("N66").Value = -20762
$ws.Range("H88").Value = 1993.1305
$ws.Range("I88").Value = 564.5714
$ws.Range("J88").Value = 2618.125
$ws.Range("K88").Value = 564.5714
$ws.Range("L88").Value = 2618.125
$ws.Range("M88").Value = -158.5714
$ws.Range("N88").Value = -3430.125
$ws.Range("H91").Value = 1993.1305
$ws.Range("I91").Value = 564.5714
$ws.Range("J91").Value = 2618.125
$ws.Range("K91").Value = 564.5714
$ws.Range("L91").Value = 2618.125
$ws.Range("M91").Value = 839.4286
$ws.Range("N91").Value = -5426.125

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 352
$ws.Range("I11").Value = 352
$ws.Range("K11").Value = 352
$ws.Range("M11").Value = -212
$ws.Range("H86").Value = 2380.889
$ws.Range("I86").Value = 1456.1538
$ws.Range("K86").Value = 1456.1538
$ws.Range("M86").Value = -333.1538
$ws.Range("H89").Value = 2380.889
$ws.Range("I89").Value = 1456.1538
$ws.Range("K89").Value = 7280.769
$ws.Range("M89").Value = -1664.769

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("M47").ClearContents()
$ws.Range("N47").ClearContents()
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()
$ws.Range("H58").Value = 2107.5908
$ws.Range("J58").Value = 2146.75
$ws.Range("L58").Value = 2146.75
$ws.Range("N58").Value = -2552.75
$ws.Range("H92").Value = 54222
$ws.Range("J92").Value = 54222
$ws.Range("L92").Value = 54222
$ws.Range("N92").Value = -59214
$ws.Range("H106").Value = 23247.25
$ws.Range("J106").Value = 23247.25
$ws.Range("L106").Value = 23247.25
$ws.Range("N106").Value = -25771.25
$ws.Range("H136").Value = 2107.5908
$ws.Range("J136").Value = 2146.75
$ws.Range("L136").Value = 6440.25
$ws.Range("N136").Value = -11540.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1953.6207
$ws.Range("I4").Value = 1967.8125
$ws.Range("J4").Value = 1936.1538
$ws.Range("K4").Value = 5903.4375
$ws.Range("L4").Value = 5808.4614
$ws.Range("M4").Value = -5791.4375
$ws.Range("N4").Value = -6032.4614
$ws.Range("H13").Value = 300
$ws.Range("I13").Value = 300
$ws.Range("K13").Value = 900
$ws.Range("M13").Value = -732
$ws.Range("H38").Value = 25.333334
$ws.Range("I38").Value = 13.666667
$ws.Range("J38").Value = 37
$ws.Range("K38").Value = 41.000001
$ws.Range("L38").Value = 111
$ws.Range("M38").Value = 305.999999
$ws.Range("N38").Value = -805

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4269.385
$ws.Range("I7").Value = 3001.75
$ws.Range("J7").Value = 4832.778
$ws.Range("K7").Value = 3001.75
$ws.Range("L7").Value = 4832.778
$ws.Range("M7").Value = -2889.75
$ws.Range("N7").Value = -5056.778
$ws.Range("H46").Value = 2949.5
$ws.Range("I46").Value = 2949.5
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 2949.5
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -2761.5
$ws.Range("N46").ClearContents()
$ws.Range("H74").Value = 24750
$ws.Range("I74").Value = 24750
$ws.Range("K74").Value = 24750
$ws.Range("M74").Value = -23752
$ws.Range("H77").Value = 24750
$ws.Range("I77").Value = 24750
$ws.Range("K77").Value = 74250
$ws.Range("M77").Value = -69258
$ws.Range("H126").Value = 4269.385
$ws.Range("I126").Value = 3001.75
$ws.Range("J126").Value = 4832.778
$ws.Range("K126").Value = 9005.25
$ws.Range("L126").Value = 14498.334
$ws.Range("M126").Value = -6535.25
$ws.Range("N126").Value = -19438.334
